$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Fonts")
$ws.Activate()

$ws.Rows("18:18").Insert()

$ws.Range("A18").Value = "group.team.point-total"
$ws.Range("B18").Value = "TradeGothicLTStd-Light.otf"
$ws.Range("C18").Value = "calibri.ttf"
$ws.Range("I18").Value = "YuGothL_0.ttf"
$ws.Range("J18").Value = "calibri.ttf"

$tbl = $ws.ListObjects.Item("Table1")
$tbl.Resize($ws.Range("A1:J25"))

$ws.Range("B18").Select()
